$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right above the old row 166, pushing the
# existing rows 166:249 down to 168:251 (matches the diff: the former
# row 166 now lives at 168, former 249 now lives at 251, etc.).
$ws.Rows.Item(166).Insert()
$ws.Rows.Item(167).Insert()

# New row 166: "Primera" quality entry for date 44510 (2021-11-10)
$ws.Range("A166").Value = 3
$ws.Range("B166").Value = "Femacal de La Calera"
$ws.Range("C166").Value = "Coquimbo"
$ws.Range("D166").Value = 44510
$ws.Range("E166").Value = 5
$ws.Range("F166").Value = "Fruta"
$ws.Range("G166").Value = 100108
$ws.Range("H166").Value = "Tropicales y subtropicales"
$ws.Range("I166").Value = 100108002
$ws.Range("J166").Value = "Mango"
$ws.Range("K166").Value = "Sin especificar"
$ws.Range("L166").Value = "Primera"
$ws.Range("M166").Value = 456
$ws.Range("N166").Value = 7000
$ws.Range("O166").Value = 7000
$ws.Range("P166").Value = 7000
$ws.Range("Q166").Value = "`$/bandeja 4 kilos"
$ws.Range("R166").Value = "Perú"
$ws.Range("S166").Value = 1750
$ws.Range("T166").Value = 4

# New row 167: "Segunda" quality entry for the same date 44510
$ws.Range("A167").Value = 3
$ws.Range("B167").Value = "Femacal de La Calera"
$ws.Range("C167").Value = "Coquimbo"
$ws.Range("D167").Value = 44510
$ws.Range("E167").Value = 5
$ws.Range("F167").Value = "Fruta"
$ws.Range("G167").Value = 100108
$ws.Range("H167").Value = "Tropicales y subtropicales"
$ws.Range("I167").Value = 100108002
$ws.Range("J167").Value = "Mango"
$ws.Range("K167").Value = "Sin especificar"
$ws.Range("L167").Value = "Segunda"
$ws.Range("M167").Value = 456
$ws.Range("N167").Value = 7000
$ws.Range("O167").Value = 7000
$ws.Range("P167").Value = 7000
$ws.Range("Q167").Value = "`$/bandeja 4 kilos"
$ws.Range("R167").Value = "Perú"
$ws.Range("S167").Value = 1750
$ws.Range("T167").Value = 4
